$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""

$ws.Range("H23").Value = 19000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""

$ws.Range("H74").Value = 16120.704
$ws.Range("I74").Value = 16263.047
$ws.Range("K74").Value = 16263.047
$ws.Range("M74").Value = -15327.047

$ws.Range("H77").Value = 16120.704
$ws.Range("I77").Value = 16263.047
$ws.Range("K77").Value = 81315.235
$ws.Range("M77").Value = -76635.235

$ws.Range("H96").Value = 1763.6
$ws.Range("I96").Value = 857
$ws.Range("J96").Value = 2556.875
$ws.Range("K96").Value = 2571
$ws.Range("L96").Value = 7670.625
$ws.Range("M96").Value = -1198
$ws.Range("N96").Value = -10416.625

$ws.Range("H107").Value = 2001.6666
$ws.Range("I107").Value = 2001.6666
$ws.Range("K107").Value = 2001.6666
$ws.Range("M107").Value = -81.66660000000002

$ws.Range("H138").Value = 2783.3386
$ws.Range("I138").Value = 2146.7932
$ws.Range("J138").Value = 3342.7273
$ws.Range("K138").Value = 6440.3796
$ws.Range("L138").Value = 10028.1819
$ws.Range("M138").Value = -1300.3796
$ws.Range("N138").Value = -20308.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5134.32
$ws.Range("I32").Value = 3459.743
$ws.Range("J32").Value = 9041.667
$ws.Range("K32").Value = 3459.743
$ws.Range("L32").Value = 9041.667
$ws.Range("M32").Value = -3172.743
$ws.Range("N32").Value = -9615.667

$ws.Range("H63").Value = 1924.2307
$ws.Range("I63").Value = 1563
$ws.Range("K63").Value = 1563
$ws.Range("M63").Value = -877

$ws.Range("H66").Value = 1924.2307
$ws.Range("I66").Value = 1563
$ws.Range("K66").Value = 7815
$ws.Range("M66").Value = -4383

$ws.Range("H132").Value = 1890350
$ws.Range("I132").Value = 2086973.6
$ws.Range("K132").Value = 6260920.800000001
$ws.Range("M132").Value = -6258390.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9999
$ws.Range("J20").Value = 9999
$ws.Range("L20").Value = 9999
$ws.Range("N20").Value = -10493

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

$ws.Range("H94").Value = 657.1667
$ws.Range("I94").Value = 641.2
$ws.Range("K94").Value = 641.2
$ws.Range("M94").Value = -190.2

$ws.Range("H105").Value = 3999.3333
$ws.Range("I105").Value = 3499.5
$ws.Range("J105").Value = 4999
$ws.Range("K105").Value = 3499.5
$ws.Range("L105").Value = 4999
$ws.Range("M105").Value = -1752.5
$ws.Range("N105").Value = -8493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 198900.67
$ws.Range("J112").Value = 198900.67
$ws.Range("L112").Value = 198900.67
$ws.Range("N112").Value = -201854.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 155.05882
$ws.Range("J2").Value = 230.125
$ws.Range("L2").Value = 1380.75
$ws.Range("N2").Value = -1606.75

$ws.Range("H103").Value = 457.7143
$ws.Range("J103").Value = 457.7143
$ws.Range("L103").Value = 1373.1429
$ws.Range("N103").Value = -3131.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 28333.334
$ws.Range("I29").Value = 50000
$ws.Range("J29").Value = 17500
$ws.Range("K29").Value = 50000
$ws.Range("L29").Value = 17500
$ws.Range("M29").Value = -49710
$ws.Range("N29").Value = -18080

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""

$ws.Range("H80").Value = 1534.25
$ws.Range("I80").Value = 1193.1
$ws.Range("J80").Value = 2102.8333
$ws.Range("K80").Value = 1193.1
$ws.Range("L80").Value = 2102.8333
$ws.Range("M80").Value = -195.0999999999999
$ws.Range("N80").Value = -4098.8333

$ws.Range("H83").Value = 1534.25
$ws.Range("I83").Value = 1193.1
$ws.Range("J83").Value = 2102.8333
$ws.Range("K83").Value = 5965.5
$ws.Range("L83").Value = 10514.1665
$ws.Range("M83").Value = -973.5
$ws.Range("N83").Value = -20498.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1641.375
$ws.Range("J22").Value = 2266.5
$ws.Range("L22").Value = 2266.5
$ws.Range("N22").Value = -2856.5

$ws.Range("H27").Value = 1641.375
$ws.Range("J27").Value = 2266.5
$ws.Range("L27").Value = 2266.5
$ws.Range("N27").Value = -2480.5

$ws.Range("H46").Value = 1239.7142
$ws.Range("I46").Value = 1112.8334
$ws.Range("K46").Value = 1112.8334
$ws.Range("M46").Value = -924.8334

$ws.Range("H82").Value = 1007.8333
$ws.Range("I82").Value = 1250.2858
$ws.Range("J82").Value = 668.4
$ws.Range("K82").Value = 1250.2858
$ws.Range("L82").Value = 668.4
$ws.Range("M82").Value = -889.2858000000001
$ws.Range("N82").Value = -1390.4

$ws.Range("H85").Value = 1007.8333
$ws.Range("I85").Value = 1250.2858
$ws.Range("J85").Value = 668.4
$ws.Range("K85").Value = 1250.2858
$ws.Range("L85").Value = 668.4
$ws.Range("M85").Value = -2.285800000000108
$ws.Range("N85").Value = -3164.4

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

$ws.Range("H100").Value = 11670888
$ws.Range("I100").Value = 17502832
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 17502832
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -17502291
$ws.Range("N100").Value = -8082

$ws.Range("H110").Value = 99999
$ws.Range("J110").Value = 99999
$ws.Range("L110").Value = 99999
$ws.Range("N110").Value = -108179

$ws.Range("H122").Value = 4213.3076
$ws.Range("I122").Value = 4325
$ws.Range("K122").Value = 12975
$ws.Range("M122").Value = -10525

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 129323
$ws.Range("J87").Value = 129323
$ws.Range("L87").Value = 129323
$ws.Range("N87").Value = -131819

$ws.Range("H90").Value = 129323
$ws.Range("J90").Value = 129323
$ws.Range("L90").Value = 387969
$ws.Range("N90").Value = -400449

$ws.Range("H112").Value = 39499.5
$ws.Range("J112").Value = 39499.5
$ws.Range("L112").Value = 39499.5
$ws.Range("N112").Value = -42453.5
